$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.869.24'
$ws.Range('E2').Value = '  -1.02%  '

$ws.Range('D3').Value = '2.464.83'
$ws.Range('E3').Value = '  -1.18%  '

$ws.Range('E4').Value = '  -0.02%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '557.80'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -2.21%  '

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '162.17'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -2.80%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('E8').Value = '  -1.19%  '

$ws.Range('D9').Value = '2.464.74'
$ws.Range('E9').Value = '  -1.15%  '

$ws.Range('E10').Value = '  -6.44%  '

$ws.Range('E11').Value = '  -0.74%  '

$ws.Range('E12').Value = '  -5.09%  '

$ws.Range('E13').Value = '  -1.39%  '

$ws.Range('D14').Value = '2.914.58'
$ws.Range('E14').Value = '  -1.25%  '

$ws.Range('D15').Value = '68.752.18'
$ws.Range('E15').Value = '  -0.95%  '

$ws.Range('E16').Value = '  -3.65%  '

$ws.Range('E17').Value = '  -3.21%  '

$ws.Range('D18').Value = '2.463.08'
$ws.Range('E18').Value = '  -0.84%  '

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '10.74'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -4.61%  '

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '341.76'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -3.73%  '

$ws.Range('E21').Value = '  -4.95%  '

$ws.Range('E22').Value = '  -2.95%  '

$ws.Range('E24').Value = '  +0.46%  '

$ws.Range('E25').Value = '  -2.65%  '

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '66.83'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -3.82%  '

$ws.Range('E27').Value = '  -4.28%  '

$ws.Range('D28').Value = '2.590.76'
$ws.Range('E28').Value = '  -1.25%  '

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +0.55%  '

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '8.12'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -5.86%  '

$ws.Range('E31').Value = '  -6.51%  '

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '7.17'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -6.06%  '

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '437.12'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -0.43%  '

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +0.00%  '

$ws.Range('E35').Value = '  -4.92%  '

$ws.Range('E36').Value = '  -6.03%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '157.25'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +2.07%  '

$ws.Range('E38').Value = '  -0.14%  '

$ws.Range('E39').Value = '  +0.13%  '

$ws.Range('E40').Value = '  -3.20%  '

$ws.Range('E42').Value = '  -3.59%  '

$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '4.44'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -3.61%  '

$ws.Range('E44').Value = '  -0.91%  '

$ws.Range('E45').Value = '  -6.79%  '

$ws.Range('E46').Value = '  +2.94%  '

$ws.Range('E47').Value = '  -5.67%  '

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '132.67'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -4.34%  '

$ws.Range('E49').Value = '  -2.69%  '

$ws.Range('E50').Value = '  -1.28%  '

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.483'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -4.64%  '
